$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Feb 08 16:41:21 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 16:41:32 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 16:41:44 EST 2023"
